$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete sales rows (rows 3-12); only the header and the
# first data row (row 2) remain after this edit.
$ws.Range("A3:A12").EntireRow.Delete()

# Row 2 updates: new timestamp, "No registrado" customer, recalculated total.
$ws.Range("B2").Value = 45777.89310185185

$ws.Range("C2").Value = "No registrado"

# E2 holds its total as text (e.g. "3200000.00"), so force a Text format
# before assigning, otherwise Excel would auto-convert the numeric-looking
# string into a number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "19250.00"
$ws.Range("E2").NumberFormat = "General"
